# Listas sem duplicação de professores
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "-"
$ws.Range("B4").Value = "-"
$ws.Range("B6").Value = "-"
$ws.Range("B7").Value = "-"

$ws.Range("C18").Value = "-"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "[-, 'ELM-2NA-Sistemas de Refrigeração', -, -]"

$ws.Range("D19").Value = "-"
$ws.Range("F19").Value = "['ELM-2NA-Sistemas de Refrigeração', -, -, -]"

$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "-"

$ws.Range("B21").Value = "-"
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "[-, -, -, 'ELM-2NA-Sistemas de Refrigeração']"
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "[-, 'ELM-2NA-Sistemas de Refrigeração', -, -]"
